$wb = $excel.ActiveWorkbook

# --- Step 1: start a new weekly sprint sheet "2017.12.12" by copying last
# week's sheet ("2017.05.12") as a template, right after it ---
$prevWeek = $wb.Worksheets.Item("2017.05.12")
$prevWeek.Copy($null, $prevWeek)
$newWeek = $wb.Worksheets.Item($prevWeek.Index + 1)
$newWeek.Name = "2017.12.12"

# First task on the new sheet: the tutoring request work item
$newWeek.Range("B3").Value2 = "Tutoring request"

# --- Step 2: close out last week's sheet - both tasks got fully worked ---
$prevWeek.Range("B3").Value2 = "Release and planning"
$prevWeek.Range("E3").Value2 = 6

# Second task on the new sheet: images for tutors
$newWeek.Range("B8").Value2 = "Images for tutors"

$prevWeek.Range("E8").Value2 = 4
$prevWeek.Range("E9").Value2 = 2

# Leave the selection on last week's sheet where the edit happened
$prevWeek.Range("B8").Select()

# Leave the new sheet's selection where work will continue
$newWeek.Range("G25").Select()

# --- Step 3: log the new sprint in the Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A21").Copy($summary.Range("A22"))
$summary.Range("A22").Value2 = 43081
$summary.Range("B21").Copy($summary.Range("B22"))
$summary.Range("B22").Value2 = "Tutoring request, images"

$summary.Activate()
$summary.Range("B22").Select()
